$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = "Dr. Shimaa Ahmad Mekki, Dr. Amira Sobhy, Dr. Rana Abo-Zaid, Dr. Hend Mahmoud, Dr. Alshimaa Atef, Dr. Heba Mahmoud Ali, Dr. Veronia Rafat, Dr. Servinaz Sayed Mohammad"
$ws.Range("G3").Value = "Dr. Alshimaa Atef, Dr. Gehan Adel, Administrator, Dr. Manar Montaser"
$ws.Range("G4").Value = "Dr. Shimaa Ahmad Mekki, Dr. Menna tuâ€™Allah Medhat, Dr. Asmaa Reda, Dr. Hanan Ragab, Dr. Majorelle Magdy, Dr. Heba Mahmoud Ali, Dr. Nourhan Mahmoud"
$ws.Range("G9").Value = "Dr. Eman M. Abo-Sakaya, Dr. Madeha Saeed, Dr. Yasmeena Fattoh, Dr. Marina Youhanna"
$ws.Range("G10").Value = "Dr. Amira Ibrahim, Dr. Basma Hamed"
$ws.Range("G12").Value = "Dr. Mona Ibrahim Hussein, Dr. Heba Al-Sayed Mohammad, Dr. Dalia Tarek Elwan"
$ws.Range("G17").Value = "Dr. Walaa Ghanima, Dr. Enas Omran"
$ws.Range("G18").Value = "Dr. Ola Abd Al-Fattah, Dr. Wafaa Ebida, Dr. Eman Samir Gabry, Dr. Abdullah El-Agrody"
$ws.Range("G19").Value = "Dr. Yasmin, Dr. Neveen Nashaat, Dr. Marina Sorial, Dr. Wafaa Ebida, Dr. Eman Samir Gabry"
$ws.Range("G20").Value = "Dr. Yasmin, Dr. Nardine, Dr. Neveen Nashaat, Dr. Marina Atef, Dr. Marina Sorial, Dr. Remon, Dr. Monica"
$ws.Range("G21").Value = "Dr. Shimaa Ahmad Mekki, Dr. Amira Sobhy, Dr. Rana Abo-Zaid, Dr. Hend Mahmoud, Dr. Alshimaa Atef, Dr. Heba Mahmoud Ali, Dr. Veronia Rafat, Dr. Servinaz Sayed Mohammad"
$ws.Range("G22").Value = "Dr. Alshimaa Atef, Dr. Gehan Adel, Administrator, Dr. Manar Montaser"
$ws.Range("G23").Value = "Dr. Shimaa Ahmad Mekki, Dr. Menna tuâ€™Allah Medhat, Dr. Asmaa Reda, Dr. Hanan Ragab, Dr. Majorelle Magdy, Dr. Heba Mahmoud Ali, Dr. Nourhan Mahmoud"
$ws.Range("G24").Value = "Dr. Abeer Ragab, Dr. Amera Ahmad Saad, Dr. Lamiaa Ossama, Dr. Nada Mohammad, Dr. Fatma Elhady"
$ws.Range("G28").Value = "Dr. Esraa Mostafa, Dr. Marwa Mustafa, Dr. Sarah Abdelmohsen, Dr. Arwa Al-Sayed, Dr. Nourhan Osama, Dr. Dina Adel, Dr. Madeha Saeed, Dr. Basma Hamed, Dr. Yasmeena Fattoh, Dr. Eman M. Abo-Sakaya"
$ws.Range("G29").Value = "Dr. Amira Ibrahim, Dr. Esraa Mostafa, Dr. Yasmeena Fattoh"
$ws.Range("G31").Value = "Dr. Mona Ibrahim Hussein, Dr. Heba Al-Sayed Mohammad, Dr. Dalia Tarek Elwan"
$ws.Range("G32").Value = "Dr. Nouran Mahmoud, Menna tuâ€™Allah Gamil"
$ws.Range("G36").Value = "Dr. Walaa Ghanima, Dr. Enas Omran"
$ws.Range("G37").Value = "Dr. Ola Abd Al-Fattah, Dr. Wafaa Ebida, Dr. Eman Samir Gabry, Dr. Abdullah El-Agrody"
$ws.Range("G38").Value = "Dr. Yasmin, Dr. Nardine, Dr. Neveen Nashaat, Dr. Marina Atef, Dr. Marina Sorial, Dr. Remon, Dr. Monica"
$ws.Range("G39").Value = "Dr. Yasmin, Dr. Nardine, Dr. Neveen Nashaat, Dr. Marina Atef, Dr. Marina Sorial, Dr. Remon, Dr. Monica"
$ws.Range("G40").Value = "Dr. Shimaa Ahmad Mekki, Dr. Amira Sobhy, Dr. Rana Abo-Zaid, Dr. Hend Mahmoud, Dr. Alshimaa Atef, Dr. Heba Mahmoud Ali, Dr. Veronia Rafat, Dr. Servinaz Sayed Mohammad"
$ws.Range("G41").Value = "Dr. Mohammad El-Tanany, Dr. Shimaa Ahmad Mekki, Dr. Hanan Ragab, Dr. Amira Sobhy, Dr. Hend Mahmoud, Dr. Alshimaa Atef"
$ws.Range("G42").Value = "Dr. Eman Tantawi, Dr. Shimaa Ahmad Mekki, Dr. Menna tuâ€™Allah Medhat, Dr. Alshimaa Atef, Dr. Servinaz Sayed Mohammad"
$ws.Range("G43").Value = "Dr. Abeer Ragab, Dr. Lamiaa Ossama, Dr. Amera Ahmad Saad, Dr. Nada Mohammad, Dr. Fatma Elhady, Dr. Kerelos Zareef, Dr. Menna tu'Alllah Mohammad"
$ws.Range("G47").Value = "Dr. Esraa Mostafa, Dr. Maryam Ahmad, Dr. Merna Said, Dr. Arwa Al-Sayed, Dr. Nourhan Osama, Dr. Amira Ibrahim"
$ws.Range("G48").Value = "Dr. Sarah Abdelmohsen, Dr. Merna Said, Dr. Maryam Ahmad, Dr. Amany Raafat, Dr. Yasmeena Fattoh, Dr. Eman M. Abo-Sakaya"
$ws.Range("G49").Value = "Dr. Mohammad Safwat, Dr. Mariam Toma Gerges"
$ws.Range("G50").Value = "Dr. Mona Ibrahim Hussein, Dr. Heba Al-Sayed Mohammad, Dr. Dalia Tarek Elwan"
$ws.Range("G51").Value = "Dr. Nouran Mahmoud, Menna tuâ€™Allah Gamil"
$ws.Range("G56").Value = "Dr. Ola Abd Al-Fattah, Dr. Wafaa Ebida, Dr. Eman Samir Gabry, Dr. Abdullah El-Agrody"
$ws.Range("G57").Value = "Dr. Yasmin, Dr. Nardine, Dr. Neveen Nashaat, Dr. Marina Atef, Dr. Marina Sorial, Dr. Remon, Dr. Monica"
$ws.Range("G58").Value = "Dr. Yasmin, Dr. Nardine, Dr. Neveen Nashaat, Dr. Marina Atef, Dr. Marina Sorial, Dr. Remon, Dr. Monica"
$ws.Range("G59").Value = "Dr. Mohammad El-Tanany, Dr. Nesma, Dr. Asmaa Reda, Dr. Amira Sobhy, Dr. Heba Mahmoud Ali, Dr. Nourhan Mahmoud, Dr. Servinaz Sayed Mohammad"
$ws.Range("G60").Value = "Dr. Mohammad El-Tanany, Dr. Shimaa Ahmad Mekki, Dr. Hanan Ragab, Dr. Amira Sobhy, Dr. Hend Mahmoud, Dr. Alshimaa Atef"
$ws.Range("G61").Value = "Dr. Shimaa Ahmad Mekki, Dr. Amira Sobhy, Dr. Asmaa Reda, Dr. Majorelle Magdy, Dr. Nahla Nagiub"
$ws.Range("G63").Value = "Dr. Safa Hany, Dr. Aya Saeed, Dr. Amal Awwad"
$ws.Range("G66").Value = "Dr. Dina Adel, Dr. Marina Youhanna, Dr. Amira Ibrahim, Dr. Madeha Saeed, Dr. Yasmeena Fattoh, Dr. Eman M. Abo-Sakaya"
$ws.Range("G67").Value = "Dr. Amira Ibrahim, Dr. Esraa Mostafa, Dr. Yasmeena Fattoh"
$ws.Range("G71").Value = "Dr. Nouran Mahmoud, Dr. Sarah Mahdy"
$ws.Range("G72").Value = "Dr. Nouran Mahmoud, Dr. Sarah Mahdy"
$ws.Range("G75").Value = "Dr. Ola Abd Al-Fattah, Dr. Wafaa Ebida, Dr. Eman Samir Gabry, Dr. Abdullah El-Agrody"
$ws.Range("G76").Value = "Dr. Yasmin, Dr. Neveen Nashaat, Dr. Marina Sorial, Dr. Wafaa Ebida, Dr. Eman Samir Gabry"
$ws.Range("G77").Value = "Dr. Yasmin, Dr. Nardine, Dr. Neveen Nashaat, Dr. Marina Atef, Dr. Marina Sorial, Dr. Remon, Dr. Monica"
$ws.Range("G78").Value = "Dr. Mohammad El-Tanany, Dr. Nesma, Dr. Asmaa Reda, Dr. Amira Sobhy, Dr. Heba Mahmoud Ali, Dr. Nourhan Mahmoud, Dr. Servinaz Sayed Mohammad"
$ws.Range("G79").Value = "Dr. Mohammad El-Tanany, Dr. Shimaa Ahmad Mekki, Dr. Hanan Ragab, Dr. Amira Sobhy, Dr. Hend Mahmoud, Dr. Alshimaa Atef"
$ws.Range("G80").Value = "Dr. Shimaa Ahmad Mekki, Dr. Amira Sobhy, Dr. Asmaa Reda, Dr. Majorelle Magdy, Dr. Nahla Nagiub"
$ws.Range("G81").Value = "Dr. Abeer Ragab, Dr. Amera Ahmad Saad, Dr. Lamiaa Ossama, Dr. Nada Mohammad, Dr. Fatma Elhady"
$ws.Range("G83").Value = "Dr. Safa Hany, Dr. Aya Saeed, Dr. Amal Awwad"
$ws.Range("G85").Value = "Dr. Dina Adel, Dr. Marina Youhanna, Dr. Amira Ibrahim, Dr. Madeha Saeed, Dr. Yasmeena Fattoh, Dr. Eman M. Abo-Sakaya"
$ws.Range("G86").Value = "Dr. Sarah Abdelmohsen, Dr. Merna Said, Dr. Maryam Ahmad, Dr. Amany Raafat, Dr. Yasmeena Fattoh, Dr. Eman M. Abo-Sakaya"
$ws.Range("G88").Value = "Dr. Mona Ibrahim Hussein, Dr. Heba Al-Sayed Mohammad, Dr. Dalia Tarek Elwan"
$ws.Range("G90").Value = "Dr. Nouran Mahmoud, Dr. Sarah Mahdy"
$ws.Range("G91").Value = "Dr. Nouran Mahmoud, Dr. Sarah Mahdy"
$ws.Range("G94").Value = "Dr. Ola Abd Al-Fattah, Dr. Wafaa Ebida, Dr. Eman Samir Gabry, Dr. Abdullah El-Agrody"
$ws.Range("G95").Value = "Dr. Yasmin, Dr. Neveen Nashaat, Dr. Marina Sorial, Dr. Wafaa Ebida, Dr. Eman Samir Gabry"
$ws.Range("G96").Value = "Dr. Yasmin, Dr. Nardine, Dr. Neveen Nashaat, Dr. Marina Atef, Dr. Marina Sorial, Dr. Remon, Dr. Monica"
$ws.Range("G97").Value = "Dr. Mohammad El-Tanany, Dr. Nesma, Dr. Asmaa Reda, Dr. Amira Sobhy, Dr. Heba Mahmoud Ali, Dr. Nourhan Mahmoud, Dr. Servinaz Sayed Mohammad"
$ws.Range("G98").Value = "Dr. Mohammad El-Tanany, Dr. Shimaa Ahmad Mekki, Dr. Hanan Ragab, Dr. Amira Sobhy, Dr. Hend Mahmoud, Dr. Alshimaa Atef"
$ws.Range("G99").Value = "Dr. Eman Tantawi, Dr. Shimaa Ahmad Mekki, Dr. Menna tuâ€™Allah Medhat, Dr. Alshimaa Atef, Dr. Servinaz Sayed Mohammad"
$ws.Range("G100").Value = "Dr. Abeer Ragab, Dr. Lamiaa Ossama, Dr. Amera Ahmad Saad, Dr. Nada Mohammad, Dr. Fatma Elhady, Dr. Kerelos Zareef, Dr. Menna tu'Alllah Mohammad"
$ws.Range("G101").Value = "Dr. Safa Hany, Dr. Aya Saeed, Dr. Amal Awwad"
$ws.Range("G104").Value = "Dr. Esraa Mostafa, Dr. Maryam Ahmad, Dr. Merna Said, Dr. Arwa Al-Sayed, Dr. Nourhan Osama, Dr. Amira Ibrahim"
$ws.Range("G105").Value = "Dr. Amira Ibrahim, Dr. Basma Hamed"
$ws.Range("G113").Value = "Dr. Ola Abd Al-Fattah, Dr. Wafaa Ebida, Dr. Eman Samir Gabry, Dr. Abdullah El-Agrody"
$ws.Range("G115").Value = "Dr. Yasmin, Dr. Nardine, Dr. Neveen Nashaat, Dr. Marina Atef, Dr. Marina Sorial, Dr. Remon, Dr. Monica"
